# Add new fields to the SERVICIO entity / form:
#   - TIPO DE DESARROLLO (column N)
#   - ESCENARIO OPI (column O)
#   - FLUJO (column P)
# These are inserted as 3 new columns before the current column N
# ("COD. REDMINE"), shifting everything from N onward to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns starting at column N (14), pushing existing
# columns N..BC to O..BF.
$insertRange = $ws.Range("N1:P1").EntireColumn
$insertRange.Insert()

# Fill in the headers for the three new columns. Shared-string table
# entries are created in first-use order, so write O1/P1/N1 in that
# sequence to land "ESCENARIO OPI", "FLUJO", "TIPO DE DESARROLLO" at
# shared-string indexes 55/56/57 respectively.
$ws.Range("O1").Value = "ESCENARIO OPI"
$ws.Range("P1").Value = "FLUJO"
$ws.Range("N1").Value = "TIPO DE DESARROLLO"

# Match the header style used by the rest of row 1 (bold white text on
# blue fill, thin border) by copying the style from the neighboring
# header cell.
$ws.Range("N1:P1").Style = $ws.Range("M1").Style

# New columns should inherit the column width of the column immediately
# to their left (standard Excel "insert column" behavior).
$ws.Range("N1:P1").EntireColumn.ColumnWidth = $ws.Range("M1").EntireColumn.ColumnWidth

# Update the selection / view to match the saved state.
$ws.Range("N12").Select()
$excel.ActiveWindow.ScrollColumn = 8
